# Fill in the new "Holter" resource block (rows 23-31) on the "External"
# sheet, mirroring the existing AEKG block (rows 14-22), and move the
# active selection to I31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startCols = @("F", "H", "J", "L", "N")
$endCols   = @("G", "I", "K", "M", "O")

# Row -> (A value, start time, end time). Row 28 is the "PAUSE" row and is
# handled separately below.
$rowData = @{
    23 = @(19, 0.375,               0.39583333333333331)
    24 = @(20, 0.39583333333333331, 0.41666666666666669)
    25 = @(21, 0.41666666666666669, 0.4375)
    26 = @(22, 0.4375,              0.46875)
    27 = @(23, 0.46875,             0.5)
    29 = @(25, 0.54166666666666663, 0.57291666666666663)
    30 = @(26, 0.57291666666666663, 0.60416666666666663)
    31 = @(27, 0.60416666666666663, 0.63541666666666663)
}

foreach ($r in 23..31) {

    # Clear any stale border-only formatting so the cells fall back to the
    # plain/general style, matching the target rows (only row 31 keeps its
    # pre-existing bordered style).
    if ($r -ne 31) {
        $ws.Range("A$r`:O$r").ClearFormats()
        $ws.Range("A$r`:O$r").ClearContents()
    } else {
        $ws.Range("B$r").ClearFormats()
        $ws.Range("B$r").ClearContents()
    }

    $ws.Range("B$r").Value = "Holter"
    $ws.Range("C$r").Value = 1
    $ws.Range("D$r").Value = "All"

    if ($r -eq 28) {
        # "PAUSE" row, same layout as rows 10 and 19.
        $ws.Range("A$r").Value = 24
        foreach ($col in $startCols) {
            $ws.Range("$col$r").Value = "PAUSE"
        }
        foreach ($col in $endCols) {
            $ws.Range("$col$r").NumberFormat = "h:mm:ss"
            $ws.Range("$col$r").Value = "PAUSE"
        }
    } else {
        $data = $rowData[$r]
        $ws.Range("A$r").Value = $data[0]
        foreach ($col in $startCols) {
            $ws.Range("$col$r").NumberFormat = "h:mm:ss"
            $ws.Range("$col$r").Value = $data[1]
        }
        foreach ($col in $endCols) {
            $ws.Range("$col$r").NumberFormat = "h:mm:ss"
            $ws.Range("$col$r").Value = $data[2]
        }
    }
}

# Row 31 keeps its original bordered style (s="4"/"5"), so it's left alone
# above other than clearing B31's leftover border-only style.

$ws.Range("I31").Select()
